$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.468.91'
$ws.Range("E2").Value = '  -3.68%  '
$ws.Range("D3").Value = '1.774.54'
$ws.Range("E3").Value = '  -2.68%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("D6").Value = "'306.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.79%  '
$ws.Range("D7").Value = "'0.4297"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.17%  '
$ws.Range("D8").Value = "'0.3663"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.43%  '
$ws.Range("D9").Value = "'0.07246"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.54%  '
$ws.Range("D10").Value = "'0.8483"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.39%  '
$ws.Range("D11").Value = "'20.34"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.32%  '
$ws.Range("D12").Value = '1.780.47'
$ws.Range("E12").Value = '  -3.39%  '
$ws.Range("D13").Value = "'5.265"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.46%  '
$ws.Range("D14").Value = "'6.442"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.53%  '
$ws.Range("D15").Value = "'0.06827"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.55%  '
$ws.Range("E16").Value = '  +0.20%  '
$ws.Range("D17").Value = "'79.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.01%  '
$ws.Range("D18").Value = "'0.000008698"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.34%  '
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").Value = "'15.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.88%  '
$ws.Range("D21").Value = '26.457.81'
$ws.Range("E21").Value = '  -3.95%  '
$ws.Range("D22").Value = "'5.102"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("D23").Value = "'11.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.45%  '
$ws.Range("D24").Value = '2.018.35'
$ws.Range("E24").Value = '  -2.04%  '
$ws.Range("D25").Value = "'152.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.58%  '
$ws.Range("D26").Value = "'1.852"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.01%  '
$ws.Range("D27").Value = "'18.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.76%  '
$ws.Range("E28").Value = '  -1.35%  '
$ws.Range("D29").Value = "'114.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.49%  '
$ws.Range("D30").Value = "'1.707"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.92%  '
$ws.Range("D31").Value = "'0.08951"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.07%  '
$ws.Range("D32").Value = "'0.7269"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.87%  '
$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").Value = "'1.120"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = "'4.346"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.23%  '
$ws.Range("D35").Value = "'2.754"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.32%  '
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("D37").Value = "'1.080"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.87%  '
$ws.Range("D38").Value = "'0.05162"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.36%  '
$ws.Range("D39").Value = "'0.01896"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.44%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = "'0.4931"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.63%  '
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = "'0.1613"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.16%  '
$ws.Range("D42").Value = "'2.524"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -9.76%  '
$ws.Range("D43").Value = "'6.226"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.44%  '
$ws.Range("D44").Value = "'8.066"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.27%  '
$ws.Range("D45").Value = "'104.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.13%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = "'1.002"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = "'10.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.34%  '
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").Value = "'0.4507"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.70%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = "'0.06206"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.17%  '
$ws.Range("D50").Value = "'1.584"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.76%  '
$ws.Range("D51").Value = "'1.745"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.88%  '